# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 586
$ws1.Range("F3").Value = 198
$ws1.Range("F4").Value = 428
$ws1.Range("F5").Value = 468
$ws1.Range("F6").Value = 276
$ws1.Range("F7").Value = 2504
$ws1.Range("F9").Value = 6657
$ws1.Range("F11").Value = 427

# Sheet "全部类型": row -> new F value (same events, different row layout)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 586
$ws4.Range("F3").Value = 198
$ws4.Range("F4").Value = 428
$ws4.Range("F5").Value = 468
$ws4.Range("F6").Value = 276
$ws4.Range("F9").Value = 2504
$ws4.Range("F11").Value = 6657
$ws4.Range("F13").Value = 427
